$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast error table values for rows 7-11 (A values 11-15)
# Row 7 (A7 = 11)
$ws.Range("B7").Value = 0.05307768655401681
$ws.Range("C7").Value = 1.531245982582867
$ws.Range("D7").Value = 5.066149181679623
$ws.Range("E7").Value = 2.250810783180057
$ws.Range("F7").Value = 2.280390043409765
$ws.Range("G7").Value = 38

# Row 8 (A8 = 12)
$ws.Range("B8").Value = 0.1368748258908087
$ws.Range("C8").Value = 1.689563907413856
$ws.Range("D8").Value = 6.751460318728611
$ws.Range("E8").Value = 2.598357234625103
$ws.Range("F8").Value = 2.630540962681023
$ws.Range("G8").Value = 37

# Row 9 (A9 = 13)
$ws.Range("B9").Value = 0.04527589459500063
$ws.Range("C9").Value = 1.985187473423592
$ws.Range("D9").Value = 9.297717034348192
$ws.Range("E9").Value = 3.04921580645716
$ws.Range("F9").Value = 3.128084520280375
$ws.Range("G9").Value = 20

# Row 10 (A10 = 14)
$ws.Range("B10").Value = -0.1230781156919481
$ws.Range("C10").Value = 2.056295129145246
$ws.Range("D10").Value = 9.159562667202785
$ws.Range("E10").Value = 3.026476939810179
$ws.Range("F10").Value = 3.147451188135231
$ws.Range("G10").Value = 13

# Row 11 (A11 = 15)
$ws.Range("B11").Value = -0.4218489719809924
$ws.Range("C11").Value = 2.349783436184865
$ws.Range("D11").Value = 10.70900605577121
$ws.Range("E11").Value = 3.272461773003805
$ws.Range("F11").Value = 3.628196780187402
$ws.Range("G11").Value = 5

$wb.Save()
